$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 142.21428
$ws.Range("I55").Value = 150.1
$ws.Range("J55").Value = 122.5
$ws.Range("K55").Value = 150.1
$ws.Range("L55").Value = 122.5
$ws.Range("M55").Value = 63.90000000000001
$ws.Range("N55").Value = -550.5

# Row 86
$ws.Range("H86").Value = 10157.632
$ws.Range("I86").Value = 10499.75
$ws.Range("J86").Value = 8333
$ws.Range("K86").Value = 10499.75
$ws.Range("L86").Value = 8333
$ws.Range("M86").Value = -9376.75
$ws.Range("N86").Value = -10579

# Row 89
$ws.Range("H89").Value = 10157.632
$ws.Range("I89").Value = 10499.75
$ws.Range("J89").Value = 8333
$ws.Range("K89").Value = 52498.75
$ws.Range("L89").Value = 41665
$ws.Range("M89").Value = -46882.75
$ws.Range("N89").Value = -52897

# Row 92
$ws.Range("H92").Value = 560.04
$ws.Range("I92").Value = 543.3
$ws.Range("K92").Value = 543.3
$ws.Range("M92").Value = 704.7

# Row 100
$ws.Range("H100").Value = 6000
$ws.Range("I100").Value = 5999.3335
$ws.Range("J100").Value = 6006
$ws.Range("K100").Value = 5999.3335
$ws.Range("L100").Value = 6006
$ws.Range("M100").Value = -5458.3335
$ws.Range("N100").Value = -7088

# Row 138
$ws.Range("H138").Value = 1829.7667
$ws.Range("J138").Value = 3163.5454
$ws.Range("L138").Value = 9490.636200000001
$ws.Range("N138").Value = -19770.6362

$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 498
$ws.Range("J6").Value = 498
$ws.Range("L6").Value = 498
$ws.Range("N6").Value = -844

# Row 11
$ws.Range("H11").Value = 83337.336
$ws.Range("J11").Value = 83337.336
$ws.Range("L11").Value = 83337.336
$ws.Range("N11").Value = -83625.336

# Row 13
$ws.Range("H13").Value = 43001.332
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 51001.6
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 51001.6
$ws.Range("M13").Value = -2856
$ws.Range("N13").Value = -51289.6

# Row 110
$ws.Range("H110").Value = 1963.6666
$ws.Range("I110").Value = 1842.6923
$ws.Range("K110").Value = 1842.6923
$ws.Range("M110").Value = 202.3077000000001

# Row 132
$ws.Range("H132").Value = 5381
$ws.Range("I132").Value = 4041.25
$ws.Range("J132").Value = 6571.8887
$ws.Range("K132").Value = 12123.75
$ws.Range("L132").Value = 19715.6661
$ws.Range("M132").Value = -9593.75
$ws.Range("N132").Value = -24775.6661

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 5819.3335
$ws.Range("I5").Value = 8504.5
$ws.Range("J5").Value = 449
$ws.Range("K5").Value = 8504.5
$ws.Range("L5").Value = 449
$ws.Range("M5").Value = -8391.5
$ws.Range("N5").Value = -675

# Row 109
$ws.Range("H109").Value = 79989.664
$ws.Range("J109").Value = 79989.664
$ws.Range("L109").Value = 79989.664
$ws.Range("N109").Value = -82763.664

# Row 134
$ws.Range("H134").Value = 29569.29
$ws.Range("I134").Value = 1409.1154
$ws.Range("J134").Value = 90583
$ws.Range("K134").Value = 4227.3462
$ws.Range("L134").Value = 271749
$ws.Range("M134").Value = -1692.3462
$ws.Range("N134").Value = -276819

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 813845.2
$ws.Range("I31").Value = 18914.777
$ws.Range("J31").Value = 1171563.9
$ws.Range("K31").Value = 18914.777
$ws.Range("L31").Value = 1171563.9
$ws.Range("M31").Value = -18619.777
$ws.Range("N31").Value = -1172153.9

# Row 34
$ws.Range("H34").Value = 813845.2
$ws.Range("I34").Value = 18914.777
$ws.Range("J34").Value = 1171563.9
$ws.Range("K34").Value = 18914.777
$ws.Range("L34").Value = 1171563.9
$ws.Range("M34").Value = -18712.777
$ws.Range("N34").Value = -1171967.9

# Row 92
$ws.Range("H92").Value = 42433.332
$ws.Range("J92").Value = 42433.332
$ws.Range("L92").Value = 42433.332
$ws.Range("N92").Value = -47425.332

# Row 107
$ws.Range("H107").Value = 777.7646999999999
$ws.Range("I107").Value = 777.7646999999999
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 777.7646999999999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1142.2353
$ws.Range("N107").ClearContents()

# Row 122
$ws.Range("H122").Value = 3395.111
$ws.Range("I122").Value = 3337
$ws.Range("J122").Value = 3598.5
$ws.Range("K122").Value = 10011
$ws.Range("L122").Value = 10795.5
$ws.Range("M122").Value = -7561
$ws.Range("N122").Value = -15695.5

# Row 132
$ws.Range("H132").Value = 2953.1177
$ws.Range("I132").Value = 2517.4167
$ws.Range("K132").Value = 7552.250100000001
$ws.Range("M132").Value = -5022.250100000001

# Row 134
$ws.Range("H134").Value = 504395.5
$ws.Range("I134").Value = 770523.9
$ws.Range("J134").Value = 10157.143
$ws.Range("K134").Value = 2311571.7
$ws.Range("L134").Value = 30471.429
$ws.Range("M134").Value = -2309036.7
$ws.Range("N134").Value = -35541.429

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 5875.2
$ws.Range("I3").Value = 2680.889
$ws.Range("J3").Value = 10666.667
$ws.Range("K3").Value = 8042.667
$ws.Range("L3").Value = 32000.001
$ws.Range("M3").Value = -7930.667
$ws.Range("N3").Value = -32224.001

# Row 37
$ws.Range("H37").Value = 82998.664
$ws.Range("J37").Value = 82998.664
$ws.Range("L37").Value = 248995.992
$ws.Range("N37").Value = -249219.992

# Row 51
$ws.Range("H51").Value = 17606.54
$ws.Range("I51").Value = 8934.375
$ws.Range("K51").Value = 26803.125
$ws.Range("M51").Value = -26343.125

# Row 80
$ws.Range("H80").Value = 4184.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4184.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 12554.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -14426.25

# Row 83
$ws.Range("H83").Value = 4184.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4184.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 37662.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -47022.75

# Row 113
$ws.Range("H113").Value = 1246.8462
$ws.Range("J113").Value = 1591.2858
$ws.Range("L113").Value = 4773.857400000001
$ws.Range("N113").Value = -9113.857400000001

# Row 122
$ws.Range("H122").Value = 1405.8
$ws.Range("I122").Value = 980
$ws.Range("J122").Value = 1547.7333
$ws.Range("K122").Value = 8820
$ws.Range("L122").Value = 13929.5997
$ws.Range("M122").Value = -6370
$ws.Range("N122").Value = -18829.5997

# Row 137
$ws.Range("H137").Value = 5792.273
$ws.Range("I137").Value = 8206
$ws.Range("K137").Value = 24618
$ws.Range("M137").Value = -19518

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 235
$ws.Range("I2").Value = 47.583332
$ws.Range("K2").Value = 47.583332
$ws.Range("M2").Value = 65.416668

# Row 3
$ws.Range("H3").Value = 3224
$ws.Range("I3").Value = 10300
$ws.Range("J3").Value = 865.3333
$ws.Range("K3").Value = 10300
$ws.Range("L3").Value = 865.3333
$ws.Range("M3").Value = -10184
$ws.Range("N3").Value = -1097.3333

# Row 12
$ws.Range("H12").Value = 174499.17
$ws.Range("J12").Value = 13665
$ws.Range("L12").Value = 13665
$ws.Range("N12").Value = -13945

# Row 70
$ws.Range("H70").Value = 4900
$ws.Range("I70").Value = 4900
$ws.Range("K70").Value = 4900
$ws.Range("M70").Value = -4630

# Row 73
$ws.Range("H73").Value = 4900
$ws.Range("I73").Value = 4900
$ws.Range("K73").Value = 4900
$ws.Range("M73").Value = -3964

# Row 122
$ws.Range("H122").Value = 1970.3043
$ws.Range("I122").Value = 1968.4286
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 5905.2858
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -3455.2858
$ws.Range("N122").Value = -10870

# Row 132
$ws.Range("H132").Value = 71432370
$ws.Range("I132").Value = 76926980
$ws.Range("K132").Value = 230780940
$ws.Range("M132").Value = -230778410

$ws = $wb.Worksheets.Item("LTW")
# Row 12
$ws.Range("H12").Value = 47429.715
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 66001.60000000001
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 66001.60000000001
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -66341.60000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 2749.5
$ws.Range("J10").Value = 4999
$ws.Range("L10").Value = 4999
$ws.Range("N10").Value = -5337

# Row 12
$ws.Range("H12").Value = 12666.333
$ws.Range("I12").Value = 16500
$ws.Range("K12").Value = 16500
$ws.Range("M12").Value = -16358

# Row 13
$ws.Range("H13").Value = 1999
$ws.Range("J13").Value = 1999
$ws.Range("L13").Value = 1999
$ws.Range("N13").Value = -2279

# Row 122
$ws.Range("H122").Value = 5682.32
$ws.Range("I122").Value = 2445.1428
$ws.Range("K122").Value = 7335.428400000001
$ws.Range("M122").Value = -4885.428400000001
